$wb = $excel.ActiveWorkbook

# Map of sheet name -> row number -> new value for column F ("想去人数")
$updates = @{
    "展览" = @{
        2  = 14172
        5  = 1221
        6  = 1056
        7  = 13966
        8  = 15123
        18 = 26
        21 = 1173
        24 = 5873
        25 = 951
        27 = 5479
        28 = 66
        30 = 74
        31 = 372
    }
    "全部类型" = @{
        2  = 14172
        6  = 1221
        7  = 1056
        8  = 13966
        9  = 15123
        19 = 26
        22 = 1173
        26 = 5873
        27 = 951
        29 = 5479
        30 = 66
        32 = 74
        33 = 372
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $updates[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $ws.Cells.Item([int]$row, 6).Value = $rowsMap[$row]
    }
}
